$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for all existing data rows (2-223)
# from 45205 to 45206.
for ($r = 2; $r -le 223; $r++) {
    $ws.Cells.Item($r, 3).Value2 = 45206
}

# Row 223 gains an explicit row height (matches default 15, but becomes customHeight).
$ws.Rows.Item(223).RowHeight = 15

# Copy formatting for the new row 224 from row 223 so styles (date format on B/C,
# wrap-text on R) match exactly.
$ws.Cells.Item(223, 2).Copy($ws.Cells.Item(224, 2))
$ws.Cells.Item(223, 3).Copy($ws.Cells.Item(224, 3))
$ws.Cells.Item(223, 18).Copy($ws.Cells.Item(224, 18))

# New row 224 data.
$ws.Cells.Item(224, 1).Value = "A 48110-2023"
$ws.Cells.Item(224, 2).Value2 = 45205
$ws.Cells.Item(224, 3).Value2 = 45206
$ws.Cells.Item(224, 4).Value = "VÄSTRA GÖTALANDS LÄN"
$ws.Cells.Item(224, 5).Value = "HERRLJUNGA"
$ws.Cells.Item(224, 7).Value = 13.2
$ws.Cells.Item(224, 8).Value = 0
$ws.Cells.Item(224, 9).Value = 0
$ws.Cells.Item(224, 10).Value = 0
$ws.Cells.Item(224, 11).Value = 0
$ws.Cells.Item(224, 12).Value = 0
$ws.Cells.Item(224, 13).Value = 0
$ws.Cells.Item(224, 14).Value = 0
$ws.Cells.Item(224, 15).Value = 0
$ws.Cells.Item(224, 16).Value = 0
$ws.Cells.Item(224, 17).Value = 0
